$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new cell value to B12 (new shared string)
$ws.Range("B12").Value = "hollllaaaaaaaaaaaaaaaa"

# Update the selection to D20 (matches the diff's <selection activeCell="D20" sqref="D20"/>)
$ws.Range("D20").Select()
